$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "Lama1"
$row2[0,2] = "Itga7"
$row2[0,3] = "ECs"
$row2[0,4] = 1
$row2[0,5] = 0.3333333333333333
$row2[0,6] = 0.014112
$row2[0,7] = 0.042336
$row2[0,8] = 0.1773673913134555
$row2[0,9] = 0.1773673913134555
$row2[0,10] = 3
$row2[0,11] = 1
$row2[0,12] = 4.851622
$row2[0,13] = 14.554866
$row2[0,14] = 0.1304267807868642
$row2[0,15] = 0.1304267807868642
$row2[0,16] = 0.068466089664
$row2[0,17] = 0.616194806976
$row2[0,18] = 0.02313345786557802
$row2[0,19] = 0.02313345786557802
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = "ECs"
$row3[0,1] = "Lama1"
$row3[0,2] = "Itga7"
$row3[0,3] = "FAPs"
$row3[0,4] = 1
$row3[0,5] = 0.3333333333333333
$row3[0,6] = 0.014112
$row3[0,7] = 0.042336
$row3[0,8] = 0.1773673913134555
$row3[0,9] = 0.1773673913134555
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 2.189198666666667
$row3[0,13] = 6.567596
$row3[0,14] = 0.05885251047922296
$row3[0,15] = 0.05885251047922296
$row3[0,16] = 0.030893971584
$row3[0,17] = 0.278045744256
$row3[0,18] = 0.01043851625594758
$row3[0,19] = 0.01043851625594758
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = "ECs"
$row4[0,1] = "Lama1"
$row4[0,2] = "Itga7"
$row4[0,3] = "Inflammatory-Mac"
$row4[0,4] = 1
$row4[0,5] = 0.3333333333333333
$row4[0,6] = 0.014112
$row4[0,7] = 0.042336
$row4[0,8] = 0.1773673913134555
$row4[0,9] = 0.1773673913134555
$row4[0,10] = 3
$row4[0,11] = 1
$row4[0,12] = 0.1764276666666667
$row4[0,13] = 0.5292830000000001
$row4[0,14] = 0.004742927747683409
$row4[0,15] = 0.004742927747683409
$row4[0,16] = 0.002489747232
$row4[0,17] = 0.022407725088
$row4[0,18] = 0.0008412407217948094
$row4[0,19] = 0.0008412407217948092
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,20
$row5[0,0] = "ECs"
$row5[0,1] = "Lama1"
$row5[0,2] = "Itga7"
$row5[0,3] = "MuSCs"
$row5[0,4] = 1
$row5[0,5] = 0.3333333333333333
$row5[0,6] = 0.014112
$row5[0,7] = 0.042336
$row5[0,8] = 0.1773673913134555
$row5[0,9] = 0.1773673913134555
$row5[0,10] = 3
$row5[0,11] = 1
$row5[0,12] = 28.76236666666667
$row5[0,13] = 86.28710000000001
$row5[0,14] = 0.7732224176048222
$row5[0,15] = 0.7732224176048221
$row5[0,16] = 0.4058945184
$row5[0,17] = 3.6530506656
$row5[0,18] = 0.1371444431156506
$row5[0,19] = 0.1371444431156506
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,20
$row6[0,0] = "ECs"
$row6[0,1] = "Lama1"
$row6[0,2] = "Itga7"
$row6[0,3] = "Neutrophils"
$row6[0,4] = 1
$row6[0,5] = 0.3333333333333333
$row6[0,6] = 0.014112
$row6[0,7] = 0.042336
$row6[0,8] = 0.1773673913134555
$row6[0,9] = 0.1773673913134555
$row6[0,10] = 3
$row6[0,11] = 1
$row6[0,12] = 0.8418256666666667
$row6[0,13] = 2.525477
$row6[0,14] = 0.02263090811425316
$row6[0,15] = 0.02263090811425315
$row6[0,16] = 0.011879843808
$row6[0,17] = 0.106918594272
$row6[0,18] = 0.004013985135279595
$row6[0,19] = 0.004013985135279593
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,20
$row7[0,0] = "ECs"
$row7[0,1] = "Lama1"
$row7[0,2] = "Itga7"
$row7[0,3] = "Resolving-Mac"
$row7[0,4] = 1
$row7[0,5] = 0.3333333333333333
$row7[0,6] = 0.014112
$row7[0,7] = 0.042336
$row7[0,8] = 0.1773673913134555
$row7[0,9] = 0.1773673913134555
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 0.3766100000000001
$row7[0,13] = 1.12983
$row7[0,14] = 0.01012445526715414
$row7[0,15] = 0.01012445526715414
$row7[0,16] = 0.005314720320000001
$row7[0,17] = 0.04783248288000001
$row7[0,18] = 0.001795748219204905
$row7[0,19] = 0.001795748219204904
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,20
$row8[0,0] = "FAPs"
$row8[0,1] = "Lama1"
$row8[0,2] = "Itga7"
$row8[0,3] = "ECs"
$row8[0,4] = 2
$row8[0,5] = 0.6666666666666666
$row8[0,6] = 0.059263
$row8[0,7] = 0.177789
$row8[0,8] = 0.7448500362393221
$row8[0,9] = 0.7448500362393219
$row8[0,10] = 3
$row8[0,11] = 1
$row8[0,12] = 4.851622
$row8[0,13] = 14.554866
$row8[0,14] = 0.1304267807868642
$row8[0,15] = 0.1304267807868642
$row8[0,16] = 0.287521674586
$row8[0,17] = 2.587695071274
$row8[0,18] = 0.0971483923956739
$row8[0,19] = 0.09714839239567387
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,20
$row9[0,0] = "FAPs"
$row9[0,1] = "Lama1"
$row9[0,2] = "Itga7"
$row9[0,3] = "FAPs"
$row9[0,4] = 2
$row9[0,5] = 0.6666666666666666
$row9[0,6] = 0.059263
$row9[0,7] = 0.177789
$row9[0,8] = 0.7448500362393221
$row9[0,9] = 0.7448500362393219
$row9[0,10] = 3
$row9[0,11] = 1
$row9[0,12] = 2.189198666666667
$row9[0,13] = 6.567596
$row9[0,14] = 0.05885251047922296
$row9[0,15] = 0.05885251047922296
$row9[0,16] = 0.1297384805826667
$row9[0,17] = 1.167646325244
$row9[0,18] = 0.04383629456322431
$row9[0,19] = 0.04383629456322429
$ws.Range("A9:T9").Value = $row9

$row10 = New-Object 'object[,]' 1,20
$row10[0,0] = "FAPs"
$row10[0,1] = "Lama1"
$row10[0,2] = "Itga7"
$row10[0,3] = "Inflammatory-Mac"
$row10[0,4] = 2
$row10[0,5] = 0.6666666666666666
$row10[0,6] = 0.059263
$row10[0,7] = 0.177789
$row10[0,8] = 0.7448500362393221
$row10[0,9] = 0.7448500362393219
$row10[0,10] = 3
$row10[0,11] = 1
$row10[0,12] = 0.1764276666666667
$row10[0,13] = 0.5292830000000001
$row10[0,14] = 0.004742927747683409
$row10[0,15] = 0.004742927747683409
$row10[0,16] = 0.01045563280966667
$row10[0,17] = 0.09410069528700001
$row10[0,18] = 0.003532769904742474
$row10[0,19] = 0.003532769904742473
$ws.Range("A10:T10").Value = $row10

$row11 = New-Object 'object[,]' 1,20
$row11[0,0] = "FAPs"
$row11[0,1] = "Lama1"
$row11[0,2] = "Itga7"
$row11[0,3] = "MuSCs"
$row11[0,4] = 2
$row11[0,5] = 0.6666666666666666
$row11[0,6] = 0.059263
$row11[0,7] = 0.177789
$row11[0,8] = 0.7448500362393221
$row11[0,9] = 0.7448500362393219
$row11[0,10] = 3
$row11[0,11] = 1
$row11[0,12] = 28.76236666666667
$row11[0,13] = 86.28710000000001
$row11[0,14] = 0.7732224176048222
$row11[0,15] = 0.7732224176048221
$row11[0,16] = 1.704544135766667
$row11[0,17] = 15.3408972219
$row11[0,18] = 0.5759347457740081
$row11[0,19] = 0.5759347457740078
$ws.Range("A11:T11").Value = $row11

$row12 = New-Object 'object[,]' 1,20
$row12[0,0] = "FAPs"
$row12[0,1] = "Lama1"
$row12[0,2] = "Itga7"
$row12[0,3] = "Neutrophils"
$row12[0,4] = 2
$row12[0,5] = 0.6666666666666666
$row12[0,6] = 0.059263
$row12[0,7] = 0.177789
$row12[0,8] = 0.7448500362393221
$row12[0,9] = 0.7448500362393219
$row12[0,10] = 3
$row12[0,11] = 1
$row12[0,12] = 0.8418256666666667
$row12[0,13] = 2.525477
$row12[0,14] = 0.02263090811425316
$row12[0,15] = 0.02263090811425315
$row12[0,16] = 0.04988911448366667
$row12[0,17] = 0.449002030353
$row12[0,18] = 0.01685663272903023
$row12[0,19] = 0.01685663272903022
$ws.Range("A12:T12").Value = $row12

$row13 = New-Object 'object[,]' 1,20
$row13[0,0] = "FAPs"
$row13[0,1] = "Lama1"
$row13[0,2] = "Itga7"
$row13[0,3] = "Resolving-Mac"
$row13[0,4] = 2
$row13[0,5] = 0.6666666666666666
$row13[0,6] = 0.059263
$row13[0,7] = 0.177789
$row13[0,8] = 0.7448500362393221
$row13[0,9] = 0.7448500362393219
$row13[0,10] = 3
$row13[0,11] = 1
$row13[0,12] = 0.3766100000000001
$row13[0,13] = 1.12983
$row13[0,14] = 0.01012445526715414
$row13[0,15] = 0.01012445526715414
$row13[0,16] = 0.02231903843
$row13[0,17] = 0.20087134587
$row13[0,18] = 0.00754120087264316
$row13[0,19] = 0.007541200872643157
$ws.Range("A13:T13").Value = $row13

$row14 = New-Object 'object[,]' 1,20
$row14[0,0] = "MuSCs"
$row14[0,1] = "Lama1"
$row14[0,2] = "Itga7"
$row14[0,3] = "ECs"
$row14[0,4] = 1
$row14[0,5] = 0.3333333333333333
$row14[0,6] = 0.006188666666666666
$row14[0,7] = 0.018566
$row14[0,8] = 0.07778257244722256
$row14[0,9] = 0.07778257244722254
$row14[0,10] = 3
$row14[0,11] = 1
$row14[0,12] = 4.851622
$row14[0,13] = 14.554866
$row14[0,14] = 0.1304267807868642
$row14[0,15] = 0.1304267807868642
$row14[0,16] = 0.03002507135066666
$row14[0,17] = 0.270225642156
$row14[0,18] = 0.01014493052561228
$row14[0,19] = 0.01014493052561228
$ws.Range("A14:T14").Value = $row14

$row15 = New-Object 'object[,]' 1,20
$row15[0,0] = "MuSCs"
$row15[0,1] = "Lama1"
$row15[0,2] = "Itga7"
$row15[0,3] = "FAPs"
$row15[0,4] = 1
$row15[0,5] = 0.3333333333333333
$row15[0,6] = 0.006188666666666666
$row15[0,7] = 0.018566
$row15[0,8] = 0.07778257244722256
$row15[0,9] = 0.07778257244722254
$row15[0,10] = 3
$row15[0,11] = 1
$row15[0,12] = 2.189198666666667
$row15[0,13] = 6.567596
$row15[0,14] = 0.05885251047922296
$row15[0,15] = 0.05885251047922296
$row15[0,16] = 0.01354822081511111
$row15[0,17] = 0.121933987336
$row15[0,18] = 0.004577699660051084
$row15[0,19] = 0.004577699660051083
$ws.Range("A15:T15").Value = $row15

$row16 = New-Object 'object[,]' 1,20
$row16[0,0] = "MuSCs"
$row16[0,1] = "Lama1"
$row16[0,2] = "Itga7"
$row16[0,3] = "Inflammatory-Mac"
$row16[0,4] = 1
$row16[0,5] = 0.3333333333333333
$row16[0,6] = 0.006188666666666666
$row16[0,7] = 0.018566
$row16[0,8] = 0.07778257244722256
$row16[0,9] = 0.07778257244722254
$row16[0,10] = 3
$row16[0,11] = 1
$row16[0,12] = 0.1764276666666667
$row16[0,13] = 0.5292830000000001
$row16[0,14] = 0.004742927747683409
$row16[0,15] = 0.004742927747683409
$row16[0,16] = 0.001091852019777778
$row16[0,17] = 0.009826668178
$row16[0,18] = 0.0003689171211461269
$row16[0,19] = 0.0003689171211461268
$ws.Range("A16:T16").Value = $row16

$row17 = New-Object 'object[,]' 1,20
$row17[0,0] = "MuSCs"
$row17[0,1] = "Lama1"
$row17[0,2] = "Itga7"
$row17[0,3] = "MuSCs"
$row17[0,4] = 1
$row17[0,5] = 0.3333333333333333
$row17[0,6] = 0.006188666666666666
$row17[0,7] = 0.018566
$row17[0,8] = 0.07778257244722256
$row17[0,9] = 0.07778257244722254
$row17[0,10] = 3
$row17[0,11] = 1
$row17[0,12] = 28.76236666666667
$row17[0,13] = 86.28710000000001
$row17[0,14] = 0.7732224176048222
$row17[0,15] = 0.7732224176048221
$row17[0,16] = 0.1780006998444444
$row17[0,17] = 1.6020062986
$row17[0,18] = 0.06014322871516366
$row17[0,19] = 0.06014322871516364
$ws.Range("A17:T17").Value = $row17

$row18 = New-Object 'object[,]' 1,20
$row18[0,0] = "MuSCs"
$row18[0,1] = "Lama1"
$row18[0,2] = "Itga7"
$row18[0,3] = "Neutrophils"
$row18[0,4] = 1
$row18[0,5] = 0.3333333333333333
$row18[0,6] = 0.006188666666666666
$row18[0,7] = 0.018566
$row18[0,8] = 0.07778257244722256
$row18[0,9] = 0.07778257244722254
$row18[0,10] = 3
$row18[0,11] = 1
$row18[0,12] = 0.8418256666666667
$row18[0,13] = 2.525477
$row18[0,14] = 0.02263090811425316
$row18[0,15] = 0.02263090811425315
$row18[0,16] = 0.005209778442444444
$row18[0,17] = 0.046888005982
$row18[0,18] = 0.001760290249943333
$row18[0,19] = 0.001760290249943332
$ws.Range("A18:T18").Value = $row18

$row19 = New-Object 'object[,]' 1,20
$row19[0,0] = "MuSCs"
$row19[0,1] = "Lama1"
$row19[0,2] = "Itga7"
$row19[0,3] = "Resolving-Mac"
$row19[0,4] = 1
$row19[0,5] = 0.3333333333333333
$row19[0,6] = 0.006188666666666666
$row19[0,7] = 0.018566
$row19[0,8] = 0.07778257244722256
$row19[0,9] = 0.07778257244722254
$row19[0,10] = 3
$row19[0,11] = 1
$row19[0,12] = 0.3766100000000001
$row19[0,13] = 1.12983
$row19[0,14] = 0.01012445526715414
$row19[0,15] = 0.01012445526715414
$row19[0,16] = 0.002330713753333334
$row19[0,17] = 0.02097642378
$row19[0,18] = 0.0007875061753060812
$row19[0,19] = 0.0007875061753060811
$ws.Range("A19:T19").Value = $row19
